$d = $word.ActiveDocument

# NOTE: the edit "rotates" several unrelated paragraphs' text content
# through each other (per the commit's xml diff). We target each
# paragraph by its fixed structural position (Paragraphs(n)) rather than
# by searching for old text globally, because several of the old/new
# strings collide with each other across paragraphs (text that is the
# "new" value in one paragraph is the "old" value of another), so a
# naive global Find/Replace pass would cascade incorrectly.

$brk = [char]11   # manual line break -> serializes as <w:br/> between <w:t> runs

# 1) "Objetivos" body paragraph
$d.Paragraphs(6).Range.Text = "Estudo das técnicas de caracterização de propriedades elétricas, magnéticas, térmicas e ópticas de materiais."

# 2) "Docente(s) Responsável(eis)" bullet paragraph
$d.Paragraphs(9).Range.Text = "Apresentar as técnicas experimentais de caracterização de propriedades elétricas, magnéticas, térmicas e ópticas de materiais."

# 3) "Programa resumido" body paragraph
$d.Paragraphs(11).Range.Text = "Propriedades elétricas: condutividade elétrica em metais puros, ligas metálicas e semicondutores,  e supercondutores; Efeito Hall; Lei de Ohm e dependência com a temperatura." + $brk + "Propriedades magnéticas: susceptibilidade magnética e magnetização c.c. Curvas de histerese de materiais magnéticos macios. Medidas de magnetostricção." + $brk + "Propriedades térmicas dos materiais:  expansão térmica."

# 4) "Programa" body paragraph
$d.Paragraphs(14).Range.Text = "Experimentos desenvolvidos em laboratório didático, realização de relatórios para cada experimento e de testes sobre o experimento em estudo."

# 5/6/7) "Avaliação" bullet paragraph holds three labeled runs (Método /
# Critério / Norma de recuperação) mixed with bold-label runs, so we
# replace each value in place via Find scoped to that paragraph's range
# (preserves the bold "label:" runs and their own formatting). We use
# wdReplaceOne (1, not wdReplaceAll) and work right-to-left (last run
# first) so an earlier replacement's freshly-written text can never be
# re-matched/cascaded into by a later (textually-colliding) Find call.
$biblio = "HUMMEL, R. E. Electronic Properties of Materials, Springer, 2000." + $brk + "RAYMOND A. SERWAY, CLEMENT J. MOSES, CURT A. MOYER. Modern Physics 3rd Edition,  Cengage Learning, Inc., 2005." + $brk + "SOLYMAR, L.; WALSH, D. Electrical Properties of Materials, Oxford University Press, 2009." + $brk + "NICOLA A. SPALDIN, Magnetic Materials, Fundamentals and Applications, SECOND EDITION, Cambridge University Press, 2011" + $brk + "ROBERT, P. Electrical and Magnetic Properties of Materials, Artech House, 1998." + $brk + "SPEYER, R. Thermal Analysis of Materials, CRC Press, 1993."

$p17 = $d.Paragraphs(17).Range.Duplicate
$p17.Find.Execute("Aplicação de uma prova escrita e prática dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação", $true, $false, $false, $false, $false, $true, 1, $false, $biblio, 1) | Out-Null

$p17 = $d.Paragraphs(17).Range.Duplicate
$p17.Find.Execute("Média aritmética das notas dos relatórios de cada experimento", $true, $false, $false, $false, $false, $true, 1, $false, "Aplicação de uma prova escrita e prática dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação", 1) | Out-Null

$p17 = $d.Paragraphs(17).Range.Duplicate
$p17.Find.Execute("Experimentos desenvolvidos em laboratório didático, realização de relatórios para cada experimento e de testes sobre o experimento em estudo.", $true, $false, $false, $false, $false, $true, 1, $false, "Média aritmética das notas dos relatórios de cada experimento", 1) | Out-Null

# 8) "Bibliografia" body paragraph
$d.Paragraphs(19).Range.Text = "5840726 - Cristina Bormio Nunes"
